$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the D (Price) column for the data rows so that Excel
# does not auto-convert numeric-looking strings (e.g. "1.0000") into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.183.09'
$ws.Range("E2").Value = '  +5.67%  '
$ws.Range("D3").Value = '1.782.35'
$ws.Range("E3").Value = '  +3.30%  '
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '244.14'
$ws.Range("E5").Value = '  +1.24%  '
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = '0.4917'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '0.2666'
$ws.Range("E8").Value = '  +2.33%  '
$ws.Range("D9").Value = '0.06247'
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("D10").Value = '1.775.63'
$ws.Range("E10").Value = '  +2.87%  '
$ws.Range("D11").Value = '16.46'
$ws.Range("E11").Value = '  +3.87%  '
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").Value = '0.6259'
$ws.Range("E13").Value = '  +2.47%  '
$ws.Range("D14").Value = '4.632'
$ws.Range("E14").Value = '  +3.09%  '
$ws.Range("D15").Value = '79.88'
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("D16").Value = '28.145.67'
$ws.Range("E16").Value = '  +6.24%  '
$ws.Range("D17").Value = '1.0000'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = '0.9999'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '0.000007222'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '12.05'
$ws.Range("E20").Value = '  +5.89%  '
$ws.Range("D21").Value = '2.006.34'
$ws.Range("E21").Value = '  +2.84%  '
$ws.Range("D22").Value = '4.559'
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").Value = '8.735'
$ws.Range("E23").Value = '  +2.29%  '
$ws.Range("D24").Value = '5.226'
$ws.Range("E24").Value = '  +2.75%  '
$ws.Range("D25").Value = '141.42'
$ws.Range("D26").Value = '15.77'
$ws.Range("E26").Value = '  +2.96%  '
$ws.Range("D27").Value = '1.857'
$ws.Range("E27").Value = '  +5.06%  '
$ws.Range("D28").Value = '109.09'
$ws.Range("E28").Value = '  +2.53%  '
$ws.Range("D29").Value = '1.386'
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").Value = '4.175'
$ws.Range("E30").Value = '  +6.78%  '
$ws.Range("D31").Value = '0.08240'
$ws.Range("E31").Value = '  +3.38%  '
$ws.Range("D32").Value = '3.770'
$ws.Range("E32").Value = '  +2.73%  '
$ws.Range("D33").Value = '0.04890'
$ws.Range("E33").Value = '  +9.01%  '
$ws.Range("D34").Value = '1.071'
$ws.Range("E34").Value = '  +6.97%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").Value = '0.6507'
$ws.Range("E36").Value = '  +4.28%  '
$ws.Range("D37").Value = '0.9445'
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("D38").Value = '2.586'
$ws.Range("D39").Value = '2.042'
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("D40").Value = '5.933'
$ws.Range("E40").Value = '  +6.18%  '
$ws.Range("D41").Value = '0.01547'
$ws.Range("E41").Value = '  +2.37%  '
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  +0.20%  '
$ws.Range("D43").Value = '99.70'
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("D44").Value = '0.3967'
$ws.Range("E44").Value = '  +3.09%  '
$ws.Range("D45").Value = '7.170'
$ws.Range("E45").Value = '  +4.55%  '
$ws.Range("D46").Value = '0.1201'
$ws.Range("E46").Value = '  +3.99%  '
$ws.Range("E47").Value = '  +0.95%  '
$ws.Range("D48").Value = '7.986'
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("D49").Value = '1.295'
$ws.Range("E49").Value = '  +5.52%  '
$ws.Range("D50").Value = '30.57'
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("D51").Value = '52.78'
$ws.Range("E51").Value = '  +2.35%  '

# Restore the cells formatting to the original (no explicit number format)
# while keeping the text values that were just assigned.
$ws.Range("D2:D51").ClearFormats()
